# Update on programm run output files
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progress")

# Tiny precision correction on previously-recorded completion timestamps
# (I3, I6, I9) coming from the latest program run.
$ws.Cells.Item(3, 9).Value = 45912.28908844908
$ws.Cells.Item(3, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(6, 9).Value = 45912.28908844908
$ws.Cells.Item(6, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(9, 9).Value = 45912.28908844908
$ws.Cells.Item(9, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the newest analysis-run rows (two new selection passes over the
# same three verses of the word).
$word = "ਲੇਪਨ"
$verse1 = "ਚੰਦਨ ਅਗਰ ਕਪੂਰ ਲੇਪਨ ਤਿਸੁ ਸੰਗੇ ਨਹੀ ਪ੍ਰੀਤਿ ॥"
$verse2 = "ਜਟਾ ਭਸਮ ਲੇਪਨ ਕੀਆ ਕਹਾ ਗੁਫਾ ਮਹਿ ਬਾਸੁ ॥"
$verse3 = "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥"

function Add-ProgressRow {
    param($row, $verse, $page, $selectedAt)

    $ws.Cells.Item($row, 1).Value = $word
    $ws.Cells.Item($row, 2).Value = $word
    $ws.Cells.Item($row, 4).Value = $verse
    $ws.Cells.Item($row, 5).Value = $page
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = $selectedAt
    $ws.Cells.Item($row, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 8).Value = "not started"
}

Add-ProgressRow 11 $verse1 1018 45912.30713481482
Add-ProgressRow 12 $verse2 1103 45912.30713481482
Add-ProgressRow 13 $verse3 1243 45912.30713481482

Add-ProgressRow 14 $verse1 1018 45912.31589575231
Add-ProgressRow 15 $verse2 1103 45912.31589575231
Add-ProgressRow 16 $verse3 1243 45912.31589575231
